$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the table header timestamp
$ws.Range("A1").Value = "Datos actualizados a 4 de Septiembre de 2020 a las 15:07"

# Updated country rows (values refreshed from source; a few countries swapped
# rank/position relative to their neighbour, e.g. Irak/Alemania, Paises Bajos/
# Emiratos Arabes Unidos and Gambia/Sri Lanka).
$countryRows = @(
    @{ Row=4; A="Estados Unidos"; B=6336138; C=894; D=3575823; E=2569201; F=0; G=56; H=191114 },
    @{ Row=6; A="India"; B=3948247; C=15123; D=3043588; E=835977; F=0; G=113; H=68682 },
    @{ Row=13; A="Argentina"; B=451198; C=0; D=331621; E=110109; F=0; G=107; H=9468 },
    @{ Row=18; A="Arabia Saudita"; B=319141; C=822; D=295063; E=20063; F=0; G=33; H=4015 },
    @{ Row=23; A="Irak"; B=252075; C=5036; D=191368; E=53348; F=0; G=84; H=7359 },
    @{ Row=24; A="Alemania"; B=249063; C=249; D=224600; E=15064; F=0; G=0; H=9399 },
    @{ Row=38; A="Kuwait"; B=88243; C=865; D=79417; E=8289; F=0; G=1; H=537 },
    @{ Row=42; A="Suecia"; B=84985; C=0; D=0; E=0; F=0; G=5; H=5835 },
    @{ Row=44; A="Paises Bajos"; B=73208; C=744; D=0; E=0; F=0; G=2; H=6237 },
    @{ Row=45; A="Emiratos Arabes Unidos"; B=72766; C=612; D=63158; E=9221; F=0; G=0; H=387 },
    @{ Row=46; A="Bielorrusia"; B=72485; C=183; D=71510; E=274; F=0; G=5; H=701 },
    @{ Row=51; A="Portugal"; B=59457; C=406; D=42576; E=15048; F=0; G=4; H=1833 },
    @{ Row=64; A="Uzbekistan"; B=42998; C=310; D=40392; E=2268; F=0; G=7; H=338 },
    @{ Row=69; A="Serbia"; B=31772; C=96; D=30387; E=664; F=0; G=3; H=721 },
    @{ Row=75; A="Estado de Palestina"; B=25142; C=671; D=16437; E=8535; F=0; G=3; H=170 },
    @{ Row=82; A="Dinamarca"; B=17547; C=173; D=15586; E=1334; F=0; G=1; H=627 },
    @{ Row=84; A="Libia"; B=16445; C=672; D=1910; E=14273; F=0; G=8; H=262 },
    @{ Row=86; A="Republica de Macedonia"; B=14871; C=109; D=12054; E=2208; F=0; G=3; H=609 },
    @{ Row=128; A="Gambia"; B=3120; C=19; D=1295; E=1726; F=0; G=0; H=99 },
    @{ Row=129; A="Sri Lanka"; B=3111; C=0; D=2907; E=192; F=0; G=0; H=12 },
    @{ Row=143; A="Islandia"; B=2135; C=7; D=2029; E=96; F=0; G=0; H=10 },
    @{ Row=179; A="Islas Feroe"; B=413; C=1; D=379; E=34; F=0; G=0; H=0 }
)

foreach ($r in $countryRows) {
    $ws.Range("A$($r.Row)").Value = $r.A
    $ws.Range("B$($r.Row)").Value = $r.B
    $ws.Range("C$($r.Row)").Value = $r.C
    $ws.Range("D$($r.Row)").Value = $r.D
    $ws.Range("E$($r.Row)").Value = $r.E
    $ws.Range("F$($r.Row)").Value = $r.F
    $ws.Range("G$($r.Row)").Value = $r.G
    $ws.Range("H$($r.Row)").Value = $r.H
}
